{"js": "// Replace the date line and every \"NNN\u00d7N=\" problem text with its updated\n// value, in document order. Every old value in this document is unique,\n// and none of the replacement values collide with any other old value,\n// so a simple ordered sequence of exact search/replace operations is\n// safe and unambiguous.\nconst replacements = [\n  [\"2024-10-03 Thursday\", \"2024-10-04 Friday\"],\n  [\"883\u00d75=\", \"527\u00d76=\"],\n  [\"195\u00d78=\", \"485\u00d79=\"],\n  [\"883\u00d73=\", \"658\u00d79=\"],\n  [\"435\u00d79=\", \"625\u00d79=\"],\n  [\"525\u00d72=\", \"586\u00d74=\"],\n  [\"961\u00d75=\", \"380\u00d75=\"],\n  [\"767\u00d73=\", \"509\u00d77=\"],\n  [\"872\u00d78=\", \"614\u00d76=\"],\n  [\"729\u00d78=\", \"143\u00d75=\"],\n  [\"251\u00d74=\", \"970\u00d73=\"],\n  [\"519\u00d74=\", \"500\u00d77=\"],\n  [\"547\u00d77=\", \"264\u00d76=\"],\n  [\"811\u00d74=\", \"449\u00d75=\"],\n  [\"516\u00d76=\", \"921\u00d73=\"],\n  [\"517\u00d74=\", \"118\u00d77=\"],\n  [\"519\u00d76=\", \"295\u00d72=\"],\n  [\"563\u00d74=\", \"855\u00d76=\"],\n  [\"621\u00d79=\", \"948\u00d76=\"],\n  [\"101\u00d78=\", \"456\u00d73=\"],\n  [\"224\u00d79=\", \"929\u00d77=\"],\n  [\"722\u00d74=\", \"584\u00d78=\"],\n  [\"228\u00d73=\", \"164\u00d72=\"],\n  [\"719\u00d75=\", \"393\u00d78=\"],\n  [\"122\u00d73=\", \"258\u00d78=\"],\n  [\"715\u00d73=\", \"676\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  // Every source value is unique in this document, so there should be\n  // exactly one hit; replace it (and any further exact duplicates,\n  // defensively) in place, preserving the run's formatting.\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"NNN\u00d7N=\" problem text with its updated\n# value, in document order. Every old value in this document is unique,\n# and none of the replacement values collide with any other old value,\n# so a simple ordered sequence of exact Find/Replace operations is safe\n# and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-10-03 Thursday\", \"2024-10-04 Friday\"),\n  @(\"883\u00d75=\", \"527\u00d76=\"),\n  @(\"195\u00d78=\", \"485\u00d79=\"),\n  @(\"883\u00d73=\", \"658\u00d79=\"),\n  @(\"435\u00d79=\", \"625\u00d79=\"),\n  @(\"525\u00d72=\", \"586\u00d74=\"),\n  @(\"961\u00d75=\", \"380\u00d75=\"),\n  @(\"767\u00d73=\", \"509\u00d77=\"),\n  @(\"872\u00d78=\", \"614\u00d76=\"),\n  @(\"729\u00d78=\", \"143\u00d75=\"),\n  @(\"251\u00d74=\", \"970\u00d73=\"),\n  @(\"519\u00d74=\", \"500\u00d77=\"),\n  @(\"547\u00d77=\", \"264\u00d76=\"),\n  @(\"811\u00d74=\", \"449\u00d75=\"),\n  @(\"516\u00d76=\", \"921\u00d73=\"),\n  @(\"517\u00d74=\", \"118\u00d77=\"),\n  @(\"519\u00d76=\", \"295\u00d72=\"),\n  @(\"563\u00d74=\", \"855\u00d76=\"),\n  @(\"621\u00d79=\", \"948\u00d76=\"),\n  @(\"101\u00d78=\", \"456\u00d73=\"),\n  @(\"224\u00d79=\", \"929\u00d77=\"),\n  @(\"722\u00d74=\", \"584\u00d78=\"),\n  @(\"228\u00d73=\", \"164\u00d72=\"),\n  @(\"719\u00d75=\", \"393\u00d78=\"),\n  @(\"122\u00d73=\", \"258\u00d78=\"),\n  @(\"715\u00d73=\", \"676\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace). Wrap=1 (wdFindContinue), Replace=2\n    # (wdReplaceAll) replaces every exact match in the document (there\n    # is exactly one occurrence of each old value here).\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
